$wb = $excel.ActiveWorkbook

# --- Sheet 1 ("Metadata"): bump version, date, contact ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B3").Value = "2.0.0"
$ws1.Range("B8").Value = "2024-06-04T14:59:10+02:00"
$ws1.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# --- Sheet 2 ("Include from FSIII"): insert 5 new concept rows ---
$ws2 = $wb.Worksheets.Item(2)

# Push the existing concept rows (J1, J5, J2, J3, J4, blank, System URI) down
# by 5 rows, making room at rows 2-6 for the new concepts.
$ws2.Rows("2:6").Insert()

# The newly inserted rows don't inherit the data-row style automatically;
# copy it over from the row that used to sit right below them.
$ws2.Range("A7:B7").Copy()
$ws2.Range("A2:B6").PasteSpecial(-4122)

# Fill in the four (well, five) new FBOE concepts.
$ws2.Range("A2").Value = "43c2b7f0-5e55-4627-8fcf-bdaf5a9d84ac"
$ws2.Range("A3").Value = "86b53158-6d05-412e-ad55-2e1fa26359b3"
$ws2.Range("A4").Value = "1c850a09-aa49-4fae-9354-f932f13e030b"
$ws2.Range("A5").Value = "462f9352-0129-4d8e-8c75-a6dfed78ddcf"
$ws2.Range("A6").Value = "4571f168-a92a-4caf-8dc8-35f45c2a1cb4"
